$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.449.60"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "2.090.39"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'329.62"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "'0.5199"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").Value = "'0.4360"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").Value = "'53.92"
$ws.Range("E9").Value = "  +16.45%  "
$ws.Range("E10").Value = "  -2.52%  "
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("D12").Value = "'24.14"
$ws.Range("E12").Value = "  -4.15%  "
$ws.Range("D13").Value = "2.091.47"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("D15").Value = "'7.655"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "'95.61"
$ws.Range("E16").Value = "  -2.90%  "
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("D19").Value = "'0.06586"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "'19.21"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D23").Value = "30.492.11"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").Value = "'12.19"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "'2.340"
$ws.Range("E25").Value = "  +3.60%  "
$ws.Range("D26").Value = "2.332.18"
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("D27").Value = "'22.18"
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("D28").Value = "'2.548"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'162.39"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").Value = "'131.30"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("D31").Value = "'1.182"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Value = "'1.644"
$ws.Range("E33").Value = "  +7.01%  "
$ws.Range("D34").Value = "'6.137"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("D35").Value = "'3.903"
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("D36").Value = "'10.11"
$ws.Range("E36").Value = "  +5.63%  "
$ws.Range("D37").Value = "'0.02571"
$ws.Range("D38").Value = "'0.06792"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "'5.438"
$ws.Range("E39").Value = "  -3.23%  "
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").Value = "'0.6860"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "'1.256"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6321"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'13.83"
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("D47").Value = "'2.189"
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("D48").Value = "'3.626"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").Value = "'1.235"
$ws.Range("E49").Value = "  +7.61%  "
$ws.Range("E50").Value = "  -4.26%  "
$ws.Range("D51").Value = "'81.53"
$ws.Range("E51").Value = "  -2.12%  "
